$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 to I1:J1
# before setting their values, so both new header cells reuse the
# existing "header" cell style (same as B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), rows 2-34
$data = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(8, 8)
    5  = @(8, 9)
    6  = @(5, 5)
    7  = @(10, 10)
    8  = @(6, 7)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(7, 7)
    15 = @(10, 10)
    16 = @(7, 8)
    17 = @(6, 7)
    18 = @(7, 7)
    19 = @(7, 8)
    20 = @(8, 9)
    21 = @(6, 7)
    22 = @(9, 9)
    23 = @(9, 9)
    24 = @(9, 9)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(8, 9)
    28 = @(12, 12)
    29 = @(3, 3)
    30 = @(3, 3)
    31 = @(9, 9)
    32 = @(2, 3)
    33 = @(7, 8)
    34 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
